$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 6844
$ws.Range("K3").Value = 7078
$ws.Range("K4").Value = 1461
$ws.Range("K5").Value = 502
$ws.Range("K6").Value = 7750
$ws.Range("K7").Value = 23635

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K4").Value = 86
$ws.Range("K7").Value = 718
$ws.Range("K8").Value = 1547
$ws.Range("K11").Value = 442
$ws.Range("K12").Value = 43
$ws.Range("K13").Value = 33
$ws.Range("K18").Value = 158
$ws.Range("K19").Value = 695
$ws.Range("K20").Value = 572
$ws.Range("K21").Value = 76
$ws.Range("K23").Value = 234
$ws.Range("K29").Value = 1288
$ws.Range("K30").Value = 91
$ws.Range("K32").Value = 25
$ws.Range("K33").Value = 1018
$ws.Range("K37").Value = 798
$ws.Range("K41").Value = 162
$ws.Range("K42").Value = 871
$ws.Range("K43").Value = 191
$ws.Range("K44").Value = 197
$ws.Range("K47").Value = 160
$ws.Range("K48").Value = 306
$ws.Range("K52").Value = 624
$ws.Range("K53").Value = 301
$ws.Range("K55").Value = 252
$ws.Range("K56").Value = 26
$ws.Range("K57").Value = 90
$ws.Range("K64").Value = 147
$ws.Range("K65").Value = 547
$ws.Range("K66").Value = 72
$ws.Range("K67").Value = 916
$ws.Range("K69").Value = 52
$ws.Range("K72").Value = 119
$ws.Range("K73").Value = 211
$ws.Range("K75").Value = 72
$ws.Range("K78").Value = 276
$ws.Range("K82").Value = 29
$ws.Range("K83").Value = 503
$ws.Range("K85").Value = 1089
$ws.Range("K88").Value = 257
$ws.Range("K91").Value = 279
$ws.Range("K94").Value = 316
$ws.Range("K95").Value = 391
$ws.Range("K97").Value = 183
$ws.Range("K99").Value = 404
$ws.Range("K101").Value = 23635

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K4").Value = 26
$ws.Range("K6").Value = 194
$ws.Range("K7").Value = 718

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K6").Value = 151
$ws.Range("K7").Value = 442

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 357
$ws.Range("K7").Value = 1089

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 167
$ws.Range("K7").Value = 624

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("K3").Value = 14
$ws.Range("K7").Value = 52

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K6").Value = 128
$ws.Range("K7").Value = 301

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 432
$ws.Range("K3").Value = 470
$ws.Range("K6").Value = 511
$ws.Range("K7").Value = 1547

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K6").Value = 116
$ws.Range("K7").Value = 503

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K3").Value = 359
$ws.Range("K4").Value = 52
$ws.Range("K6").Value = 322
$ws.Range("K7").Value = 1018

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K2").Value = 134
$ws.Range("K3").Value = 136
$ws.Range("K4").Value = 16
$ws.Range("K7").Value = 391

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K5").Value = 33
$ws.Range("K7").Value = 798

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 180
$ws.Range("K3").Value = 134
$ws.Range("K7").Value = 547

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K2").Value = 107
$ws.Range("K3").Value = 168
$ws.Range("K7").Value = 404

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("K2").Value = 26
$ws.Range("K6").Value = 32
$ws.Range("K7").Value = 91

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K5").Value = 22
$ws.Range("K6").Value = 264
$ws.Range("K7").Value = 916

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 365
$ws.Range("K3").Value = 458
$ws.Range("K6").Value = 374
$ws.Range("K7").Value = 1288

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K2").Value = 47
$ws.Range("K4").Value = 41
$ws.Range("K6").Value = 144
$ws.Range("K7").Value = 306

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K3").Value = 208
$ws.Range("K6").Value = 230
$ws.Range("K7").Value = 695

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K6").Value = 78
$ws.Range("K7").Value = 197

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K4").Value = 11
$ws.Range("K7").Value = 162

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 237
$ws.Range("K7").Value = 871

$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("K3").Value = 12
$ws.Range("K6").Value = 33

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K3").Value = 72
$ws.Range("K7").Value = 276

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K2").Value = 77
$ws.Range("K7").Value = 252

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K3").Value = 81
$ws.Range("K7").Value = 234

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K2").Value = 71
$ws.Range("K3").Value = 132
$ws.Range("K7").Value = 279

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("K3").Value = 20
$ws.Range("K7").Value = 76

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K3").Value = 41
$ws.Range("K7").Value = 147

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K6").Value = 156
$ws.Range("K7").Value = 572

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K3").Value = 53
$ws.Range("K6").Value = 41
$ws.Range("K7").Value = 158

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K6").Value = 144
$ws.Range("K7").Value = 316

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K3").Value = 48
$ws.Range("K4").Value = 15
$ws.Range("K6").Value = 49
$ws.Range("K7").Value = 160

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("K6").Value = 36
$ws.Range("K7").Value = 72

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K2").Value = 73
$ws.Range("K7").Value = 211

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K3").Value = 40
$ws.Range("K7").Value = 183

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K3").Value = 80
$ws.Range("K7").Value = 257

$ws = $wb.Worksheets.Item("Galewood")
$ws.Range("K3").Value = 7
$ws.Range("K7").Value = 25

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("K2").Value = 23
$ws.Range("K7").Value = 72

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("K6").Value = 39
$ws.Range("K7").Value = 90

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K2").Value = 38
$ws.Range("K6").Value = 72
$ws.Range("K7").Value = 191

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K6").Value = 56
$ws.Range("K7").Value = 119

$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("K4").Value = 5
$ws.Range("K6").Value = 29

$ws = $wb.Worksheets.Item("Magnificent Mile")
$ws.Range("K6").Value = 13
$ws.Range("K7").Value = 26

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("K2").Value = 28
$ws.Range("K7").Value = 86

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("K2").Value = 11
$ws.Range("K7").Value = 43
